$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB2").Value = -0.166400372770544
$ws.Range("AD2").Value = 1.206763604850376
$ws.Range("AG2").Value = -6.252890313084356
$ws.Range("AB3").Value = -0.295996312127685
$ws.Range("AD3").Value = -0.6618273521133181
$ws.Range("AG3").Value = -13.39297015020611
$ws.Range("AB4").Value = -0.3074432006063255
$ws.Range("AD4").Value = 0.2403117847406054
$ws.Range("AG4").Value = -6.506070841504995
$ws.Range("AB5").Value = -0.2902421723677454
$ws.Range("AD5").Value = -0.3813231726682013
$ws.Range("AG5").Value = -6.938865102429284
$ws.Range("AB6").Value = -0.3792063163659771
$ws.Range("AD6").Value = -2.311836642690182
$ws.Range("AG6").Value = -13.93964611521012
$ws.Range("AB7").Value = -0.2372094180570924
$ws.Range("AD7").Value = 0.06352804711987661
$ws.Range("AG7").Value = -9.339838335931347
$ws.Range("AB8").Value = -0.354056202531341
$ws.Range("AD8").Value = -2.613140207894537
$ws.Range("AG8").Value = -15.91972090079977
$ws.Range("AB9").Value = -0.3666067631597452
$ws.Range("AD9").Value = -2.476898309927298
$ws.Range("AG9").Value = -11.72441550934934
$ws.Range("AB10").Value = -0.359883292689782
$ws.Range("AD10").Value = -2.080484820666785
$ws.Range("AG10").Value = -13.04643243446868
$ws.Range("AB11").Value = -0.323263173420538
$ws.Range("AD11").Value = -1.873029929308267
$ws.Range("AG11").Value = -12.5950833077475
$ws.Range("AB12").Value = -0.3288605364191713
$ws.Range("AD12").Value = -1.541874255089881
$ws.Range("AG12").Value = -11.83356758091845
$ws.Range("AB13").Value = -0.3114362264009238
$ws.Range("AD13").Value = -0.8796961801374111
$ws.Range("AG13").Value = -11.29718986207859
$ws.Range("AB14").Value = -0.2936119744798988
$ws.Range("AD14").Value = -0.7163682132542064
$ws.Range("AG14").Value = -7.301373672670099
$ws.Range("AB15").Value = -0.336108560818545
$ws.Range("AD15").Value = -1.897432208837434
$ws.Range("AG15").Value = -10.98436080101969
$ws.Range("AB16").Value = -0.2825910882449856
$ws.Range("AD16").Value = -1.276173450199792
$ws.Range("AG16").Value = -10.63949971836509
$ws.Range("AB17").Value = -0.2822367851134115
$ws.Range("AD17").Value = -1.758261383858426
$ws.Range("AG17").Value = -14.33667450812506
$ws.Range("AB18").Value = -0.3738842957307086
$ws.Range("AD18").Value = -2.177232989355285
$ws.Range("AG18").Value = -13.31070224053161
$ws.Range("AB19").Value = -0.3059633354622712
$ws.Range("AD19").Value = -2.215781303349771
$ws.Range("AG19").Value = -13.29640927967048
$ws.Range("AB20").Value = -0.1714054313875508
$ws.Range("AD20").Value = 2.487862376675363
$ws.Range("AG20").Value = -8.877813287370689
$ws.Range("AB21").Value = -0.388626272628939
$ws.Range("AD21").Value = -0.4919399957999145
$ws.Range("AG21").Value = -11.21018977532791
$ws.Range("AB22").Value = -0.211315038441335
$ws.Range("AD22").Value = 2.276278044530738
$ws.Range("AG22").Value = -6.681038797009853
$ws.Range("AB23").Value = -0.3911413300651606
$ws.Range("AD23").Value = -0.4408995922747213
$ws.Range("AG23").Value = -8.890957381713731
$ws.Range("AB24").Value = -0.2914807327057374
$ws.Range("AD24").Value = -0.9493067988182872
$ws.Range("AG24").Value = -9.935056865408384
$ws.Range("AB25").Value = -0.1952355654546289
$ws.Range("AD25").Value = -1.330942454033567
$ws.Range("AG25").Value = -7.48286577415073
$ws.Range("AB26").Value = -0.3547353885848022
$ws.Range("AD26").Value = -0.07271006184924378
$ws.Range("AG26").Value = -12.68375295230208
$ws.Range("AB27").Value = -0.3453339062761586
$ws.Range("AD27").Value = -1.810393508255256
$ws.Range("AG27").Value = -11.79832312560712
$ws.Range("AB28").Value = -0.2726968453543164
$ws.Range("AD28").Value = -1.292425142566827
$ws.Range("AG28").Value = -13.87202914372516
$ws.Range("AB29").Value = -0.3557734116958147
$ws.Range("AD29").Value = -2.804973297146931
$ws.Range("AG29").Value = -12.49249571548325
$ws.Range("AB30").Value = -0.3330484021214341
$ws.Range("AD30").Value = -4.644302866457737
$ws.Range("AG30").Value = -13.90613421139488
$ws.Range("AB31").Value = -0.4280826510776105
$ws.Range("AD31").Value = -1.240277177657472
$ws.Range("AG31").Value = -11.2815723336375
$ws.Range("AB32").Value = -0.3149096965297397
$ws.Range("AD32").Value = -3.813562732736394
$ws.Range("AG32").Value = -13.12727029901161
$ws.Range("AB33").Value = -0.2600970374300035
$ws.Range("AD33").Value = -0.2745048171169281
$ws.Range("AG33").Value = -9.885881632361704
$ws.Range("AB34").Value = -0.12548761041987
$ws.Range("AD34").Value = 1.062250868309531
$ws.Range("AG34").Value = -7.788724894943198
$ws.Range("AB35").Value = -0.360747391990778
$ws.Range("AD35").Value = 0.3916166776330529
$ws.Range("AG35").Value = -10.34410608485515
$ws.Range("AB36").Value = -0.3396246821899999
$ws.Range("AD36").Value = -2.216615054129072
$ws.Range("AG36").Value = -14.80517331510241
$ws.Range("AB37").Value = -0.3061277679427306
$ws.Range("AD37").Value = 0.05280274697743949
$ws.Range("AG37").Value = -9.454942319978079
$ws.Range("AB38").Value = -0.2574029041840175
$ws.Range("AD38").Value = -2.37291360622711
$ws.Range("AG38").Value = -12.55492663258322
$ws.Range("AB39").Value = -0.3255151247602819
$ws.Range("AD39").Value = -0.7057648713668029
$ws.Range("AG39").Value = -13.47775037723206
$ws.Range("AB40").Value = -0.2787200484459515
$ws.Range("AD40").Value = 1.059623716409263
$ws.Range("AG40").Value = -6.789889273623598
$ws.Range("AB41").Value = -0.2927032072016966
$ws.Range("AD41").Value = -0.513950080483028
$ws.Range("AG41").Value = -13.41816047910564
$ws.Range("AB42").Value = -0.364362594078273
$ws.Range("AD42").Value = -1.817004594972412
$ws.Range("AG42").Value = -13.19922641494747
$ws.Range("AB43").Value = -0.3814164375157714
$ws.Range("AD43").Value = 0.7646933332672882
$ws.Range("AG43").Value = -9.331980803818851
$ws.Range("AB44").Value = -0.3325576634969708
$ws.Range("AD44").Value = -0.9323822980461995
$ws.Range("AG44").Value = -11.40974199729068
$ws.Range("AB45").Value = -0.33874534569127
$ws.Range("AD45").Value = -0.2571485131976691
$ws.Range("AG45").Value = -9.022946953986615
$ws.Range("AB46").Value = -0.3821023360168416
$ws.Range("AD46").Value = 0.2742320430889965
$ws.Range("AG46").Value = -11.35111968847128
$ws.Range("AB47").Value = -0.3565462312341346
$ws.Range("AD47").Value = -2.401087641333458
$ws.Range("AG47").Value = -13.58831089634392
$ws.Range("AB48").Value = -0.2638665020039602
$ws.Range("AD48").Value = -2.018414950201884
$ws.Range("AG48").Value = -13.22981040450057
$ws.Range("AB49").Value = -0.06033166280510745
$ws.Range("AD49").Value = 2.560831682811852
$ws.Range("AG49").Value = -6.158259903120213
$ws.Range("AB50").Value = -0.2531421791624721
$ws.Range("AD50").Value = -1.457118660769562
$ws.Range("AG50").Value = -9.557106847950354
$ws.Range("AB51").Value = -0.2598579162039645
$ws.Range("AD51").Value = -1.960033407080617
$ws.Range("AG51").Value = -10.37668813995711
$ws.Range("AB52").Value = -0.2853614256819633
$ws.Range("AD52").Value = -1.390947525765934
$ws.Range("AG52").Value = -12.2959537000667
$ws.Range("AB53").Value = -0.3338712546352632
$ws.Range("AD53").Value = 0.8411691084461106
$ws.Range("AG53").Value = -6.808987953061466
$ws.Range("AB54").Value = -0.3526427880242517
$ws.Range("AD54").Value = 0.2256053764421199
$ws.Range("AG54").Value = -8.250070543842785
$ws.Range("AB55").Value = -0.2464125002836758
$ws.Range("AD55").Value = -2.59709785208677
$ws.Range("AG55").Value = -12.27300827102258
$ws.Range("AB56").Value = -0.3962041294550909
$ws.Range("AD56").Value = -1.799844870719656
$ws.Range("AG56").Value = -11.88101780096861
$ws.Range("AB57").Value = -0.2746583752242959
$ws.Range("AD57").Value = -2.341694315064976
$ws.Range("AG57").Value = -11.68615830044116
$ws.Range("AB58").Value = -0.2597234112576138
$ws.Range("AD58").Value = 1.61259302329642
$ws.Range("AG58").Value = -6.245222071081882
$ws.Range("AB59").Value = -0.4200693917297115
$ws.Range("AD59").Value = -2.95011303390724
$ws.Range("AG59").Value = -13.81154128524457
$ws.Range("AB60").Value = -0.2888850143619949
$ws.Range("AD60").Value = -2.547396548106795
$ws.Range("AG60").Value = -14.21337252358856
